# Update "江西-漫展信息.xlsx":
#  - Remove the oldest convention entry (2024-07-14, 吉安·COMIC LIFE次元假日05),
#    which shifts every later row up by one position.
#  - Re-sequence the index column (A) back to 1..39 for the shifted rows.
#  - Bump the "想去人数" (interest count, column F) for a number of rows to
#    the freshly scraped values.
#
# The same edit applies identically to both the "展览" sheet and the
# "全部类型" sheet (they carry duplicate data in this workbook).

$wb = $excel.ActiveWorkbook

# Interest-count (column F) refresh, keyed by the *new* row number
# (i.e. after the 2024-07-14 row has been removed and rows shifted up).
$fUpdates = @{
    2  = 251
    3  = 1307
    4  = 140
    5  = 260
    6  = 219
    8  = 10
    9  = 172
    10 = 123
    11 = 4402
    12 = 6676
    14 = 52
    15 = 93
    16 = 559
    17 = 53
    18 = 4091
    19 = 456
    21 = 44
    22 = 2668
    25 = 161
    26 = 335
    27 = 345
    28 = 392
    29 = 215
    31 = 1606
    32 = 1010
    34 = 118
    35 = 75
    36 = 529
    38 = 10
    40 = 617
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # 1. Drop the first data row (row 2); everything below shifts up one row.
    $ws.Rows.Item(2).Delete()

    # 2. The index/id column (A) holds static numbers, not a formula, so the
    #    delete+shift does not renumber it on its own - put 1..39 back.
    for ($r = 2; $r -le 40; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $r - 1
    }

    # 3. Refresh the interest counts that changed since the last scrape.
    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value2 = $fUpdates[$r]
    }
}
